$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the formatting already used by the other header cells (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF)
$iValues = @(9, 8, 8, 8, 5, 6, 8, 9, 7, 9, 9)
$jValues = @(9, 8, 8, 8, 6, 7, 8, 9, 8, 9, 9)

for ($row = 2; $row -le 12; $row++) {
    $idx = $row - 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
